# Renal Abnormalities of Magnesium Metabolism.xlsx
# - Add a new "metadata" worksheet after the existing "data" sheet.
# - Populate it with one header row + one data row describing the PanelApp query.
# - Refresh the "time_taken" timestamps (column F) on the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update the per-gene query timestamps on the "data" sheet (col F, rows 2-20) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:35:31.993028"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:31.993036"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:31.993039"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:31.993042"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:31.993045"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:31.993047"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:31.993050"
$dataSheet.Range("F9").Value = "2021-10-05 14:35:31.993053"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:31.993055"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:31.993058"
$dataSheet.Range("F12").Value = "2021-10-05 14:35:31.993061"
$dataSheet.Range("F13").Value = "2021-10-05 14:35:31.993063"
$dataSheet.Range("F14").Value = "2021-10-05 14:35:31.993066"
$dataSheet.Range("F15").Value = "2021-10-05 14:35:31.993068"
$dataSheet.Range("F16").Value = "2021-10-05 14:35:31.993071"
$dataSheet.Range("F17").Value = "2021-10-05 14:35:31.993074"
$dataSheet.Range("F18").Value = "2021-10-05 14:35:31.993077"
$dataSheet.Range("F19").Value = "2021-10-05 14:35:31.993079"
$dataSheet.Range("F20").Value = "2021-10-05 14:35:31.993082"

# --- 2. Add the new "metadata" sheet right after "data" ---
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Match the "data" sheet's (openpyxl default) page margins: 0.75/0.75/1/1 in,
# 0.5/0.5 in header/footer == 54/54/72/72/36/36 points.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Reuse the "data" sheet's bold/centered/thin-bordered header style for B1:G1
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Reuse the "data" sheet's header-row index style for A2
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Renal Abnormalities of Magnesium Metabolism"
$ws.Range("C2").Value = 197
$ws.Range("E2").Value = "2021-06-14T08:33:47.873759Z"
$ws.Range("F2").Value = "2021-10-05 14:35:31.989187"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/197/?format=json"

# D2 ("0.28") must stay plain text (no trailing-style cell), so stamp it as
# Text first, then strip the style back to the sheet default via a
# formats-only paste from an unstyled cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.28"
$dataSheet.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$dataSheet.Activate()

Write-Output "done"
